$p = $ppt.ActivePresentation

# 1. Remove the bullet "Trabalhos relacionados" from the agenda slide (slide 2).
$agendaShape = $p.Slides.Item(2).Shapes.Item(1)
$agendaShape.TextFrame.TextRange.Paragraphs(1,1).Delete()

# 2. Hide slide 4.
$p.Slides.Item(4).SlideShowTransition.Hidden = -1

# 3. Delete the last four slides (19-22): "TECNOLOGIAS UTILIZADAS" and the
#    three flow-chart / screenshot slides at the end of the deck. Delete from
#    the end so indices of the slides being removed stay stable.
$p.Slides.Item(22).Delete()
$p.Slides.Item(21).Delete()
$p.Slides.Item(20).Delete()
$p.Slides.Item(19).Delete()
